# Apply the updated cryptocurrency price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold plain-text price strings (many look numeric, e.g. "0.556"
# or "5.20"); a leading apostrophe forces Excel to keep them as text instead of
# silently reinterpreting/rounding them as numbers. Values with multiple "." groups
# (e.g. "76.657.63") are never auto-converted, so they are written as-is.

$ws.Range("D2").Value = '76.657.63'
$ws.Range("E2").Value = '  +1.06%  '

$ws.Range("D3").Value = '3.028.52'
$ws.Range("E3").Value = '  +4.43%  '

$ws.Range("D5").Value = "'" + '202.13'
$ws.Range("E5").Value = '  +1.21%  '

$ws.Range("D6").Value = "'" + '631.82'
$ws.Range("E6").Value = '  +6.07%  '

$ws.Range("D8").Value = "'" + '0.556'
$ws.Range("E8").Value = '  +1.54%  '

$ws.Range("D9").Value = "'" + '0.212'
$ws.Range("E9").Value = '  +7.17%  '

$ws.Range("D10").Value = '3.026.07'
$ws.Range("E10").Value = '  +4.39%  '

$ws.Range("E11").Value = '  +3.44%  '

$ws.Range("E12").Value = '  -0.28%  '

$ws.Range("D13").Value = "'" + '5.22'
$ws.Range("E13").Value = '  +7.59%  '

$ws.Range("D14").Value = '3.586.21'
$ws.Range("E14").Value = '  +4.47%  '

$ws.Range("D15").Value = "'" + '29.46'
$ws.Range("E15").Value = '  +7.44%  '

$ws.Range("D16").Value = '76.512.41'
$ws.Range("E16").Value = '  +0.97%  '

$ws.Range("E17").Value = '  +2.73%  '

$ws.Range("D18").Value = '3.018.48'
$ws.Range("E18").Value = '  +3.64%  '

$ws.Range("E19").Value = '  +6.48%  '

$ws.Range("D20").Value = "'" + '8.86'
$ws.Range("E20").Value = '  -0.47%  '

$ws.Range("D21").Value = "'" + '378.32'
$ws.Range("E21").Value = '  +1.99%  '

$ws.Range("E22").Value = '  +1.20%  '

$ws.Range("D23").Value = "'" + '4.39'
$ws.Range("E23").Value = '  +2.91%  '

$ws.Range("D24").Value = "'" + '73.71'
$ws.Range("E24").Value = '  +3.94%  '

$ws.Range("D25").Value = '3.184.96'

$ws.Range("D26").Value = "'" + '4.42'
$ws.Range("E26").Value = '  +6.25%  '

$ws.Range("D27").Value = "'" + '0.998'
$ws.Range("E27").Value = '  -0.11%  '

$ws.Range("D28").Value = "'" + '9.99'
$ws.Range("E28").Value = '  +3.74%  '

$ws.Range("E29").Value = '  +3.91%  '

$ws.Range("D30").Value = "'" + '0.997'
$ws.Range("E30").Value = '  +0.07%  '

$ws.Range("D31").Value = "'" + '8.39'
$ws.Range("E31").Value = '  +8.95%  '

$ws.Range("E32").Value = '  +2.29%  '

$ws.Range("D33").Value = "'" + '513.42'
$ws.Range("E33").Value = '  +2.82%  '

$ws.Range("E34").Value = '  +8.78%  '

$ws.Range("D35").Value = "'" + '0.999'
$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("E36").Value = '  +3.76%  '

$ws.Range("D37").Value = "'" + '163.84'
$ws.Range("E37").Value = '  -0.35%  '

$ws.Range("D38").Value = "'" + '0.385'
$ws.Range("E38").Value = '  +12.18%  '

$ws.Range("D39").Value = "'" + '20.03'
$ws.Range("E39").Value = '  +2.01%  '

$ws.Range("D40").Value = "'" + '0.107'
$ws.Range("E40").Value = '  +6.20%  '

$ws.Range("D41").Value = "'" + '188.96'
$ws.Range("E41").Value = '  +4.51%  '

$ws.Range("D42").Value = "'" + '0.114'
$ws.Range("E42").Value = '  +1.34%  '

$ws.Range("E43").Value = '  +0.29%  '

$ws.Range("D44").Value = "'" + '5.20'
$ws.Range("E44").Value = '  +5.25%  '

$ws.Range("D45").Value = "'" + '42.47'
$ws.Range("E45").Value = '  +5.84%  '

$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").Value = "'" + '1.27'
$ws.Range("E46").Value = '  +7.39%  '

$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = "'" + '1.68'
$ws.Range("E47").Value = '  +3.13%  '

$ws.Range("D48").Value = "'" + '2.46'
$ws.Range("E48").Value = '  +6.56%  '

$ws.Range("D49").Value = "'" + '0.610'
$ws.Range("E49").Value = '  +7.24%  '

$ws.Range("D50").Value = "'" + '0.716'
$ws.Range("E50").Value = '  +9.79%  '

$ws.Range("E51").Value = '  +5.88%  '
